$wb = $excel.ActiveWorkbook

# ---------------------------------------------------------------------------
# Sheet "nhap-linhkien" (import log): replace the two data rows with a new
# shipment (GRM1555C1H271JA01D / S1M-13-F) and drop the extra duplicate rows
# that used to follow them.
# ---------------------------------------------------------------------------
$ws1 = $wb.Worksheets.Item("nhap-linhkien")

# Remove the now-unused trailing rows (4:6) before rewriting rows 2:3.
$ws1.Rows("4:6").Delete()

$ws1.Range("A2").Value = "GRM1555C1H271JA01D"
$ws1.Range("B2").Value = "Ceramic Cap 0402 270p 50V 5% NP0"
$ws1.Range("C2").Value = "sohopdong01"
$ws1.Range("D2").Value = "sanpham01"
$ws1.Range("E2").Value = "cty01"
$ws1.Range("F2").NumberFormat = "@"
$ws1.Range("F2").Value = "2021-10-10"
$ws1.Range("F2").Style = "Normal"
$ws1.Range("G2").Value = "cái"
$ws1.Range("H2").Value = 12
$ws1.Range("I2").NumberFormat = "@"
$ws1.Range("I2").Value = "21.00"
$ws1.Range("I2").Style = "Normal"
$ws1.Range("J2").NumberFormat = "@"
$ws1.Range("J2").Value = "252.0000"
$ws1.Range("J2").Style = "Normal"

$ws1.Range("A3").Value = "S1M-13-F"
$ws1.Range("B3").Value = "S1M"
$ws1.Range("C3").Value = "sohopdong02"
$ws1.Range("D3").Value = "sanpham02"
$ws1.Range("E3").Value = "cty02"
$ws1.Range("F3").NumberFormat = "@"
$ws1.Range("F3").Value = "2021-10-10"
$ws1.Range("F3").Style = "Normal"
$ws1.Range("G3").Value = "cái"
$ws1.Range("H3").Value = 22
$ws1.Range("I3").NumberFormat = "@"
$ws1.Range("I3").Value = "10.69"
$ws1.Range("I3").Style = "Normal"
$ws1.Range("J3").NumberFormat = "@"
$ws1.Range("J3").Value = "235.1800"
$ws1.Range("J3").Style = "Normal"

# ---------------------------------------------------------------------------
# Sheet "xuat-linhkien" (export log): the single export row is removed,
# leaving only the header row.
# ---------------------------------------------------------------------------
$ws2 = $wb.Worksheets.Item("xuat-linhkien")
$ws2.Rows("2:2").Delete()

# ---------------------------------------------------------------------------
# Sheet "ton-linhkien" (stock on hand): now reflects the two items from the
# new shipment above; the old rows 4:5 are dropped.
# ---------------------------------------------------------------------------
$ws3 = $wb.Worksheets.Item("ton-linhkien")
$ws3.Rows("4:5").Delete()

$ws3.Range("A2").Value = "Ceramic Cap 0402 270p 50V 5% NP0"
$ws3.Range("B2").Value = 12
$ws3.Range("C2").Value = "cái"

$ws3.Range("A3").Value = "S1M"
$ws3.Range("B3").Value = 22
$ws3.Range("C3").Value = "cái"
